$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.8245905782990299
